{"js": "// Date header + 24 division problems updated per the diff.\n// Each \"before\" text is unique in the document, so a targeted\n// search-and-replace per pair reproduces the edit exactly.\nconst replacements = [\n  [\"2025-06-05 Thursday\", \"2025-06-06 Friday\"],\n  [\"631\u00f79=\", \"350\u00f73=\"],\n  [\"915\u00f77=\", \"298\u00f78=\"],\n  [\"367\u00f73=\", \"544\u00f78=\"],\n  [\"133\u00f73=\", \"579\u00f79=\"],\n  [\"786\u00f78=\", \"998\u00f76=\"],\n  [\"866\u00f74=\", \"486\u00f75=\"],\n  [\"726\u00f72=\", \"725\u00f73=\"],\n  [\"447\u00f72=\", \"488\u00f75=\"],\n  [\"660\u00f75=\", \"551\u00f77=\"],\n  [\"592\u00f78=\", \"581\u00f74=\"],\n  [\"338\u00f79=\", \"866\u00f72=\"],\n  [\"783\u00f78=\", \"831\u00f72=\"],\n  [\"491\u00f78=\", \"656\u00f76=\"],\n  [\"868\u00f72=\", \"810\u00f76=\"],\n  [\"709\u00f73=\", \"177\u00f79=\"],\n  [\"370\u00f72=\", \"387\u00f79=\"],\n  [\"417\u00f73=\", \"680\u00f74=\"],\n  [\"893\u00f79=\", \"399\u00f75=\"],\n  [\"577\u00f78=\", \"427\u00f74=\"],\n  [\"695\u00f73=\", \"211\u00f75=\"],\n  [\"646\u00f78=\", \"971\u00f78=\"],\n  [\"746\u00f72=\", \"794\u00f79=\"],\n  [\"801\u00f72=\", \"859\u00f77=\"],\n  [\"319\u00f76=\", \"195\u00f76=\"],\n  [\"636\u00f78=\", \"391\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${before}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Date header + 24 division problems updated per the diff.\n# Each \"before\" text is unique in the document, so a targeted\n# Find/Replace per pair reproduces the edit exactly.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @('2025-06-05 Thursday', '2025-06-06 Friday'),\n  @('631\u00f79=', '350\u00f73='),\n  @('915\u00f77=', '298\u00f78='),\n  @('367\u00f73=', '544\u00f78='),\n  @('133\u00f73=', '579\u00f79='),\n  @('786\u00f78=', '998\u00f76='),\n  @('866\u00f74=', '486\u00f75='),\n  @('726\u00f72=', '725\u00f73='),\n  @('447\u00f72=', '488\u00f75='),\n  @('660\u00f75=', '551\u00f77='),\n  @('592\u00f78=', '581\u00f74='),\n  @('338\u00f79=', '866\u00f72='),\n  @('783\u00f78=', '831\u00f72='),\n  @('491\u00f78=', '656\u00f76='),\n  @('868\u00f72=', '810\u00f76='),\n  @('709\u00f73=', '177\u00f79='),\n  @('370\u00f72=', '387\u00f79='),\n  @('417\u00f73=', '680\u00f74='),\n  @('893\u00f79=', '399\u00f75='),\n  @('577\u00f78=', '427\u00f74='),\n  @('695\u00f73=', '211\u00f75='),\n  @('646\u00f78=', '971\u00f78='),\n  @('746\u00f72=', '794\u00f79='),\n  @('801\u00f72=', '859\u00f77='),\n  @('319\u00f76=', '195\u00f76='),\n  @('636\u00f78=', '391\u00f74='),\n)\n\nforeach ($pair in $replacements) {\n  $before = $pair[0]\n  $after = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $before\n  $find.Replacement.Text = $after\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.Execute([ref]$null, [ref]$true, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$true, [ref]$null, [ref]$null, [ref]$null, 2) | Out-Null\n}\n"}
